# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
#
# Numeric-looking text values (e.g. "317.55") are written with a leading
# apostrophe so Excel keeps them as TEXT (matching the sheet's original
# inlineStr cells) instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.084.23'
$ws.Range("E2").Value = '  -1.04%  '

# Row 3
$ws.Range("D3").Value = '1.794.94'
$ws.Range("E3").Value = '  +0.14%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '''317.55'
$ws.Range("E5").Value = '  +1.03%  '

# Row 6
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '  +0.04%  '

# Row 7
$ws.Range("D7").Value = '''0.5393'
$ws.Range("E7").Value = '  -1.25%  '

# Row 8
$ws.Range("D8").Value = '''0.3779'
$ws.Range("E8").Value = '  -1.12%  '

# Row 9
$ws.Range("D9").Value = '''0.07448'
$ws.Range("E9").Value = '  -1.54%  '

# Row 10
$ws.Range("D10").Value = '''41.74'
$ws.Range("E10").Value = '  -1.62%  '

# Row 11
$ws.Range("E11").Value = '  -2.28%  '

# Row 12
$ws.Range("E12").Value = '  -0.09%  '

# Row 13
$ws.Range("D13").Value = '''20.53'
$ws.Range("E13").Value = '  -2.51%  '

# Row 14
$ws.Range("D14").Value = '''6.111'

# Row 15
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.797.94'
$ws.Range("E15").Value = '  +0.17%  '

# Row 16
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''7.235'
$ws.Range("E16").Value = '  -2.00%  '

# Row 17
$ws.Range("D17").Value = '''89.04'
$ws.Range("E17").Value = '  -2.62%  '

# Row 18
$ws.Range("D18").Value = '''0.00001059'
$ws.Range("E18").Value = '  -0.63%  '

# Row 19
$ws.Range("E19").Value = '  +0.51%  '

# Row 20
$ws.Range("E20").Value = '  -0.05%  '

# Row 21
$ws.Range("E21").Value = '  -0.38%  '

# Row 22
$ws.Range("D22").Value = '''5.897'
$ws.Range("E22").Value = '  -0.87%  '

# Row 23
$ws.Range("D23").Value = '28.110.65'
$ws.Range("E23").Value = '  -0.94%  '

# Row 24
$ws.Range("D24").Value = '''11.16'
$ws.Range("E24").Value = '  -2.01%  '

# Row 25
$ws.Range("D25").Value = '''2.092'
$ws.Range("E25").Value = '  -1.28%  '

# Row 26
$ws.Range("D26").Value = '''155.33'
$ws.Range("E26").Value = '  -2.58%  '

# Row 27
$ws.Range("D27").Value = '''20.28'
$ws.Range("E27").Value = '  -1.89%  '

# Row 28
$ws.Range("D28").Value = '1.998.68'
$ws.Range("E28").Value = '  -0.12%  '

# Row 29
$ws.Range("E29").Value = '  -4.71%  '

# Row 30
$ws.Range("D30").Value = '''121.13'
$ws.Range("E30").Value = '  -1.49%  '

# Row 31
$ws.Range("D31").Value = '''1.117'
$ws.Range("E31").Value = '  -0.09%  '

# Row 32
$ws.Range("D32").Value = '''0.1063'
$ws.Range("E32").Value = '  +3.71%  '

# Row 33
$ws.Range("D33").Value = '''3.654'
$ws.Range("E33").Value = '  -1.09%  '

# Row 34
$ws.Range("D34").Value = '''5.547'
$ws.Range("E34").Value = '  -3.04%  '

# Row 35
$ws.Range("D35").Value = '''0.2254'
$ws.Range("E35").Value = '  -2.25%  '

# Row 36
$ws.Range("D36").Value = '''0.06468'
$ws.Range("E36").Value = '  +0.97%  '

# Row 37
$ws.Range("D37").Value = '''0.02290'
$ws.Range("E37").Value = '  -1.11%  '

# Row 38
$ws.Range("D38").Value = '''5.008'
$ws.Range("E38").Value = '  -3.04%  '

# Row 39
$ws.Range("D39").Value = '''8.443'
$ws.Range("E39").Value = '  -3.40%  '

# Row 40
$ws.Range("B40").Value = 'WEMIXTOKEN'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '''1.451'
$ws.Range("E40").Value = '  +4.63%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.6177'
$ws.Range("E41").Value = '  -2.98%  '

# Row 42
$ws.Range("D42").Value = '''11.10'
$ws.Range("E42").Value = '  -4.23%  '

# Row 43
$ws.Range("E43").Value = '  +1.73%  '

# Row 44
$ws.Range("D44").Value = '''1.001'
$ws.Range("E44").Value = '  +0.08%  '

# Row 45
$ws.Range("D45").Value = '''13.26'
$ws.Range("E45").Value = '  -1.94%  '

# Row 46
$ws.Range("E46").Value = '  +0.01%  '

# Row 47
$ws.Range("D47").Value = '''0.5775'
$ws.Range("E47").Value = '  -2.89%  '

# Row 48
$ws.Range("D48").Value = '''124.19'
$ws.Range("E48").Value = '  -1.38%  '

# Row 49
$ws.Range("D49").Value = '''1.189'
$ws.Range("E49").Value = '  +3.74%  '

# Row 50
$ws.Range("D50").Value = '''1.923'
$ws.Range("E50").Value = '  -3.17%  '

# Row 51
$ws.Range("D51").Value = '''0.06812'
$ws.Range("E51").Value = '  -1.32%  '
